$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 65.0933347143764
$ws.Range("F2").Value = 111.4415081739426
$ws.Range("I2").Value = 33.90676093101501
$ws.Range("J2").Value = 23.73373818397522
$ws.Range("K2").Value = 13.53249561786652
$ws.Range("B3").Value = 62.88228458109916
$ws.Range("F3").Value = 111.9140541553497
$ws.Range("I3").Value = 34.41623687744141
$ws.Range("J3").Value = 24.73254895210266
$ws.Range("K3").Value = 15.35366725921631
$ws.Range("B4").Value = 62.96438829882936
$ws.Range("F4").Value = 112.2823359966278
$ws.Range("I4").Value = 34.67719101905823
$ws.Range("J4").Value = 25.51641476154327
$ws.Range("K4").Value = 15.04389727115631
$ws.Range("B5").Value = 62.84512054764491
$ws.Range("F5").Value = 112.4855961799622
$ws.Range("I5").Value = 34.63815772533417
$ws.Range("J5").Value = 25.93167185783386
$ws.Range("K5").Value = 15.21130514144897
$ws.Range("B6").Value = 62.6110206978301
$ws.Range("F6").Value = 112.5782541036606
$ws.Range("I6").Value = 34.55391812324524
$ws.Range("J6").Value = 26.11307632923126
$ws.Range("K6").Value = 15.53505384922028
$ws.Range("B7").Value = 60.88687433668565
$ws.Range("F7").Value = 112.5843144655228
$ws.Range("I7").Value = 34.55053949356079
$ws.Range("J7").Value = 26.11169612407684
$ws.Range("K7").Value = 17.26353168487549
$ws.Range("B8").Value = 60.44338134656027
$ws.Range("F8").Value = 112.5670640468597
$ws.Range("I8").Value = 34.49979865550995
$ws.Range("J8").Value = 26.06681990623474
$ws.Range("K8").Value = 17.76005554199219
$ws.Range("B9").Value = 61.35889125206268
$ws.Range("F9").Value = 112.2715930938721
$ws.Range("I9").Value = 34.17254543304443
$ws.Range("J9").Value = 25.42989027500153
$ws.Range("K9").Value = 17.18862426280975
$ws.Range("B10").Value = 65.26545644046018
$ws.Range("F10").Value = 111.4937838315964
$ws.Range("I10").Value = 33.62367534637451
$ws.Range("J10").Value = 23.81232762336731
$ws.Range("K10").Value = 13.67571091651917
$ws.Range("B11").Value = 69.52732898053364
$ws.Range("F11").Value = 110.7118648290634
$ws.Range("I11").Value = 33.04508292675018
$ws.Range("J11").Value = 22.20508444309235
$ws.Range("K11").Value = 9.578214287757874
$ws.Range("B12").Value = 70.93259731582657
$ws.Range("F12").Value = 110.3363039493561
$ws.Range("I12").Value = 32.66748857498169
$ws.Range("J12").Value = 21.58015191555023
$ws.Range("K12").Value = 8.139475464820862
$ws.Range("B13").Value = 71.14714724622354
$ws.Range("F13").Value = 110.0406731367111
$ws.Range("I13").Value = 32.36169445514679
$ws.Range("J13").Value = 21.99669444561005
$ws.Range("K13").Value = 7.011850833892822
$ws.Range("B14").Value = 70.53810752055688
$ws.Range("F14").Value = 109.9915798902512
$ws.Range("I14").Value = 32.24607801437378
$ws.Range("J14").Value = 22.8040212392807
$ws.Range("K14").Value = 6.819340586662292
$ws.Range("B15").Value = 70.00658219742036
$ws.Range("F15").Value = 110.0725702047348
$ws.Range("I15").Value = 32.20220446586609
$ws.Range("J15").Value = 23.32365560531616
$ws.Range("K15").Value = 7.112635374069214
$ws.Range("B16").Value = 69.76635424906817
$ws.Range("F16").Value = 110.0665476322174
$ws.Range("I16").Value = 32.31143569946289
$ws.Range("J16").Value = 23.4561333656311
$ws.Range("K16").Value = 7.109052300453186
$ws.Range("B17").Value = 69.54118974798985
$ws.Range("F17").Value = 110.0550218820572
$ws.Range("I17").Value = 32.44962024688721
$ws.Range("J17").Value = 23.54602193832397
$ws.Range("K17").Value = 7.089839577674866
$ws.Range("B18").Value = 68.69197304974477
$ws.Range("F18").Value = 110.2991166114807
$ws.Range("I18").Value = 32.64772689342499
$ws.Range("J18").Value = 23.95540714263916
$ws.Range("K18").Value = 8.018659353256226
$ws.Range("B19").Value = 68.3466828652472
$ws.Range("F19").Value = 110.4454981088638
$ws.Range("I19").Value = 32.84017062187195
$ws.Range("J19").Value = 23.94261145591736
$ws.Range("K19").Value = 8.585562109947205
$ws.Range("B20").Value = 69.49666655869078
$ws.Range("F20").Value = 110.4458416700363
$ws.Range("I20").Value = 32.82030010223389
$ws.Range("J20").Value = 22.78108108043671
$ws.Range("K20").Value = 8.574025273323059
$ws.Range("B21").Value = 71.9577080555955
$ws.Range("F21").Value = 110.0497258901596
$ws.Range("I21").Value = 32.5433177947998
$ws.Range("J21").Value = 20.94118082523346
$ws.Range("K21").Value = 7.066733479499817
$ws.Range("B22").Value = 73.34722074041565
$ws.Range("F22").Value = 109.6521730422974
$ws.Range("I22").Value = 32.36775958538055
$ws.Range("J22").Value = 20.0566303730011
$ws.Range("K22").Value = 5.577640771865845
$ws.Range("B23").Value = 74.28690701414416
$ws.Range("F23").Value = 109.3401387929916
$ws.Range("I23").Value = 32.28172290325165
$ws.Range("J23").Value = 19.43602120876312
$ws.Range("K23").Value = 4.417166709899902
$ws.Range("B24").Value = 74.57640858193918
$ws.Range("F24").Value = 109.2202410697937
$ws.Range("I24").Value = 32.3069521188736
$ws.Range("J24").Value = 19.20245742797852
$ws.Range("K24").Value = 3.980450868606567
$ws.Range("B25").Value = 72.31220011532605
$ws.Range("F25").Value = 109.805584192276
$ws.Range("I25").Value = 32.90310072898865
$ws.Range("J25").Value = 20.40261697769165
$ws.Range("K25").Value = 7.112635374069214
